# Append rows 8-13 to the "Ruturaj Gaikwad" match log sheet, duplicating
# existing rows 2, 3, 7, 6, 5, 4 (in that order) as new rows at the bottom,
# matching the source scrape's refreshed JSON export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(" Dubai (DSC)", " October 25 2020", "Super Kings won by 8 wickets (with 8 balls remaining)", "Chennai Super Kings", "Royal Challengers Bangalore", "Ruturaj Gaikwad ", "65", "51", "4", "3", "127.45"),
    @(" Sharjah", " October 23 2020", "Mumbai won by 10 wickets (with 46 balls remaining)", "Chennai Super Kings", "Mumbai Indians", "Ruturaj Gaikwad ", "0", "5", "0", "0", "0.00"),
    @(" Abu Dhabi", " November 01 2020", "Super Kings won by 9 wickets (with 7 balls remaining)", "Chennai Super Kings", "Kings XI Punjab", "Ruturaj Gaikwad ", "62", "49", "6", "1", "126.53"),
    @(" Sharjah", " September 22 2020", "Royals won by 16 runs", "Chennai Super Kings", "Rajasthan Royals", "Ruturaj Gaikwad ", "0", "1", "0", "0", "0.00"),
    @(" Dubai (DSC)", " October 29 2020", "Super Kings won by 6 wickets", "Chennai Super Kings", "Kolkata Knight Riders", "Ruturaj Gaikwad ", "72", "53", "6", "2", "135.84"),
    @(" Dubai (DSC)", " September 25 2020", "Capitals won by 44 runs", "Chennai Super Kings", "Delhi Capitals", "Ruturaj Gaikwad ", "5", "10", "0", "0", "50.00")
)

$startRow = 8
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $values = $rows[$i]
    $rowRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 11))
    $rowRange.NumberFormat = "@"
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $values[$c]
    }
}
